# Update odds figures for the matches in row 5 (Chapecoense-SC vs Goias)
# and row 24 (Sp. Luqueno vs Nacional Asuncion) on the FlashScore sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 updates ---
$ws.Range("G5").Value = 2.7
$ws.Range("I5").Value = 2.88
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 3.6
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 7
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 26
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 11
$ws.Range("AK5").Value = 29
$ws.Range("AL5").Value = 26
$ws.Range("AN5").Value = 4.5
$ws.Range("AT5").Value = 2.38
$ws.Range("AW5").Value = 4.75
$ws.Range("AX5").Value = 17
$ws.Range("AY5").Value = 29
$ws.Range("AZ5").Value = 51

# --- Row 24 updates ---
$ws.Range("G24").Value = 2.9
$ws.Range("H24").Value = 3.1
$ws.Range("I24").Value = 2.5
$ws.Range("N24").Value = 7.5
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("W24").Value = 8
$ws.Range("AW24").Value = 4.5
